$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.007.62"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").Value = "1.742.35"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.76"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.01%  "
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5051"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -4.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2745"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06181"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").Value = "1.748.50"
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07256"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.6517"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.14"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.677"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.54"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.42%  "
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").Value = "26.019.55"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.89"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006873"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.49%  "
$ws.Range("D21").Value = "1.971.46"
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.464"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.715"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.366"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "135.85"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.492"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.30"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.781"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "105.71"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.890"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08182"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.649"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04656"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.658"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9967"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.805"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6091"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01625"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.926"
$ws.Range("D39").ClearFormats()
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "100.61"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.3916"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7634"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.998"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1164"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.320"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.59"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05304"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.64"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.68%  "
$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3458"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.52%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.582"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.17%  "
